$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the literal string into the cell as TEXT (matching the
    # source workbook, which stores these as inline strings) rather
    # than letting Excel auto-coerce number-looking text into a number.
    # =T("...") always evaluates to a text value; Copy + PasteSpecial
    # (xlPasteValues = -4163) then bakes that text in place of the
    # formula without touching the cells existing style/number format.
    $escapedForFormula = $text.Replace('"', '""')
    $rng = $ws.Range($cellRef)
    $rng.Formula = '=T("' + $escapedForFormula + '")'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

Set-TextValue "D2" '58.359.80'
Set-TextValue "E2" '  -4.25%  '
Set-TextValue "D3" '2.646.06'
Set-TextValue "E3" '  -1.67%  '
Set-TextValue "E4" '  +0.16%  '
Set-TextValue "D5" '521.71'
Set-TextValue "E5" '  -0.89%  '
Set-TextValue "D6" '144.60'
Set-TextValue "E6" '  -0.38%  '
Set-TextValue "E7" '  +0.25%  '
Set-TextValue "E8" '  -1.41%  '
Set-TextValue "E9" '  +2.06%  '
Set-TextValue "E10" '  -3.10%  '
Set-TextValue "E11" '  -0.61%  '
Set-TextValue "D12" '0.132'
Set-TextValue "E12" '  +1.64%  '
Set-TextValue "D13" '3.112.37'
Set-TextValue "E13" '  -1.59%  '
Set-TextValue "D14" '58.365.79'
Set-TextValue "E14" '  -4.10%  '
Set-TextValue "D15" '20.94'
Set-TextValue "E15" '  -1.79%  '
Set-TextValue "D16" '0.0000136'
Set-TextValue "E16" '  -1.25%  '
Set-TextValue "D17" '2.660.02'
Set-TextValue "E17" '  -1.09%  '
Set-TextValue "D18" '337.90'
Set-TextValue "E18" '  -2.92%  '
Set-TextValue "E19" '  -2.87%  '
Set-TextValue "D20" '10.47'
Set-TextValue "E20" '  -1.11%  '
Set-TextValue "D21" '6.30'
Set-TextValue "E21" '  -0.55%  '
Set-TextValue "E22" '  +0.31%  '
Set-TextValue "D23" '64.35'
Set-TextValue "E23" '  +0.97%  '
Set-TextValue "E24" '  +0.63%  '
Set-TextValue "E25" '  -2.09%  '
Set-TextValue "E26" '  +0.57%  '
Set-TextValue "D27" '0.0₃0797'
Set-TextValue "E27" '  -2.40%  '
Set-TextValue "D28" '7.13'
Set-TextValue "E28" '  -3.05%  '
Set-TextValue "D29" '6.68'
Set-TextValue "E29" '  -2.90%  '
Set-TextValue "E30" '  +0.10%  '
Set-TextValue "E31" '  -0.80%  '
Set-TextValue "D32" '153.18'
Set-TextValue "E32" '  +1.70%  '
Set-TextValue "D33" '18.86'
Set-TextValue "E33" '  -1.39%  '
Set-TextValue "E34" '  -2.79%  '
Set-TextValue "E35" '  -4.40%  '
Set-TextValue "D36" '0.908'
Set-TextValue "E36" '  -3.11%  '
Set-TextValue "D37" '0.863'
Set-TextValue "E37" '  -2.42%  '
Set-TextValue "D38" '36.87'
Set-TextValue "E38" '  -0.25%  '
Set-TextValue "E39" '  -4.64%  '
Set-TextValue "E40" '  -0.86%  '
Set-TextValue "E41" '  +0.35%  '
Set-TextValue "D42" '0.608'
Set-TextValue "E42" '  -0.85%  '
Set-TextValue "D43" '272.22'
Set-TextValue "E43" '  -5.64%  '
Set-TextValue "D44" '0.0972'
Set-TextValue "E44" '  -2.13%  '
Set-TextValue "D45" '19.41'
Set-TextValue "E45" '  -2.96%  '
Set-TextValue "E46" '  +1.54%  '
Set-TextValue "E47" '  -0.91%  '
Set-TextValue "D48" '2.049.87'
Set-TextValue "E48" '  -4.35%  '
Set-TextValue "E49" '  -2.88%  '
Set-TextValue "D50" '4.68'
Set-TextValue "E50" '  -4.11%  '
Set-TextValue "D51" '18.34'
Set-TextValue "E51" '  -4.18%  '
